$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The row that used to be row 13 held only "8554681 - Pedro Felipe Arce Castillo"
# in B13/C13 (no A13 label). Deleting it shifts every row below up by one,
# which is exactly what the target layout (dimension A1:C23, one row fewer)
# requires.
$ws.Rows.Item(13).Delete()

# After the shift, a handful of B/C cells need their text content replaced
# (the row labels in column A were already correct after the shift).

# Row 13 (label "Programa resumido:") -> "Semestral"
$ws.Range("B13:C13").Value = "Semestral"

# Row 15 (label "Programa:") -> "01/01/2012"
# Use Copy from B8:C8 (which already holds this exact text) instead of
# assigning the literal string, so Excel doesn't reinterpret it as a date.
$ws.Range("B8:C8").Copy($ws.Range("B15:C15"))

# Row 18 (label "Método:") -> "8554681 - Pedro Felipe Arce Castillo"
$ws.Range("B18:C18").Value = "8554681 - Pedro Felipe Arce Castillo"

# Row 19 (label "Critério:") -> "2 provas escritas"
$ws.Range("B19:C19").Value = "2 provas escritas"

# Row 20 (label "Norma de recuperação:") -> long evaluation-criteria text
$ws.Range("B20:C20").Value = "serão avaliados os conteúdos discutidos em sala e constantes da ementa do curso. A média da disciplina será a média aritmética das duas provas."

# Row 21 (label "Bibliografia:") -> "prova escrita com conteúdo de todo o semestre"
$ws.Range("B21:C21").Value = "prova escrita com conteúdo de todo o semestre"
